$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row (row 31) with the two new key/value strings and
# reuse the existing "True" shared string (already used by C27:C30) for C31.
$ws.Range("A31").Value = "res.users.group_multi_currency"
$ws.Range("B31").Value = "base.group_multi_currency"
$ws.Range("C27").Copy()
$ws.Range("C31").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Reflect the sheet's final UI/selection state as saved by the author
# (entire rows were selected, with B31 left as the active cell).
$ws.Rows.Select()
$ws.Range("B31").Activate()
